$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out row 1 entirely (becomes an empty row)
$ws.Range("A1:D1").ClearContents()

# Update row 2 with new lead info
$ws.Range("A2").Value = "Asfiya Wasim`r`n2022 Volkswagen Taos (New)"
$ws.Range("B2").Value = "GUBAGOO - CHAT LEAD"
$ws.Range("C2").Value = "10/30 12:57 AM"
